{"js": "// Update the division-problem answers in the table to the new values.\n// Each cell holds a single run of text like \"961\u00f79=106, 7\"; we locate the\n// exact old text with a case-sensitive, non-wildcard search and replace it\n// in place so the run's formatting (font/size) is preserved.\nconst replacements = [\n  [\"961\u00f79=106, 7\", \"533\u00f74=133, 1\"],\n  [\"608\u00f73=202, 2\", \"527\u00f75=105, 2\"],\n  [\"543\u00f78=67, 7\", \"843\u00f72=421, 1\"],\n  [\"604\u00f79=67, 1\", \"750\u00f77=107, 1\"],\n  [\"993\u00f73=331, 0\", \"315\u00f75=63, 0\"],\n  [\"982\u00f76=163, 4\", \"677\u00f79=75, 2\"],\n  [\"399\u00f76=66, 3\", \"844\u00f77=120, 4\"],\n  [\"858\u00f75=171, 3\", \"822\u00f77=117, 3\"],\n  [\"274\u00f79=30, 4\", \"949\u00f76=158, 1\"],\n  [\"928\u00f79=103, 1\", \"280\u00f76=46, 4\"],\n  [\"428\u00f73=142, 2\", \"816\u00f74=204, 0\"],\n  [\"126\u00f76=21, 0\", \"161\u00f77=23, 0\"],\n  [\"656\u00f78=82, 0\", \"545\u00f76=90, 5\"],\n  [\"634\u00f79=70, 4\", \"150\u00f79=16, 6\"],\n  [\"838\u00f77=119, 5\", \"314\u00f78=39, 2\"],\n  [\"942\u00f73=314, 0\", \"996\u00f73=332, 0\"],\n  [\"796\u00f73=265, 1\", \"425\u00f79=47, 2\"],\n  [\"892\u00f79=99, 1\", \"630\u00f72=315, 0\"],\n  [\"875\u00f79=97, 2\", \"663\u00f74=165, 3\"],\n  [\"961\u00f72=480, 1\", \"743\u00f75=148, 3\"],\n  [\"581\u00f79=64, 5\", \"776\u00f73=258, 2\"],\n  [\"106\u00f76=17, 4\", \"406\u00f78=50, 6\"],\n  [\"900\u00f73=300, 0\", \"850\u00f74=212, 2\"],\n  [\"747\u00f76=124, 3\", \"607\u00f73=202, 1\"],\n  [\"205\u00f77=29, 2\", \"199\u00f72=99, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division-problem answers in the table to the new values.\n# Each cell holds a single run of text like \"961\u00f79=106, 7\"; we use\n# Find/Replace (exact text, match case) so the run's formatting\n# (font/size) is preserved and only the text content changes.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"961\u00f79=106, 7\", \"533\u00f74=133, 1\"),\n    @(\"608\u00f73=202, 2\", \"527\u00f75=105, 2\"),\n    @(\"543\u00f78=67, 7\", \"843\u00f72=421, 1\"),\n    @(\"604\u00f79=67, 1\", \"750\u00f77=107, 1\"),\n    @(\"993\u00f73=331, 0\", \"315\u00f75=63, 0\"),\n    @(\"982\u00f76=163, 4\", \"677\u00f79=75, 2\"),\n    @(\"399\u00f76=66, 3\", \"844\u00f77=120, 4\"),\n    @(\"858\u00f75=171, 3\", \"822\u00f77=117, 3\"),\n    @(\"274\u00f79=30, 4\", \"949\u00f76=158, 1\"),\n    @(\"928\u00f79=103, 1\", \"280\u00f76=46, 4\"),\n    @(\"428\u00f73=142, 2\", \"816\u00f74=204, 0\"),\n    @(\"126\u00f76=21, 0\", \"161\u00f77=23, 0\"),\n    @(\"656\u00f78=82, 0\", \"545\u00f76=90, 5\"),\n    @(\"634\u00f79=70, 4\", \"150\u00f79=16, 6\"),\n    @(\"838\u00f77=119, 5\", \"314\u00f78=39, 2\"),\n    @(\"942\u00f73=314, 0\", \"996\u00f73=332, 0\"),\n    @(\"796\u00f73=265, 1\", \"425\u00f79=47, 2\"),\n    @(\"892\u00f79=99, 1\", \"630\u00f72=315, 0\"),\n    @(\"875\u00f79=97, 2\", \"663\u00f74=165, 3\"),\n    @(\"961\u00f72=480, 1\", \"743\u00f75=148, 3\"),\n    @(\"581\u00f79=64, 5\", \"776\u00f73=258, 2\"),\n    @(\"106\u00f76=17, 4\", \"406\u00f78=50, 6\"),\n    @(\"900\u00f73=300, 0\", \"850\u00f74=212, 2\"),\n    @(\"747\u00f76=124, 3\", \"607\u00f73=202, 1\"),\n    @(\"205\u00f77=29, 2\", \"199\u00f72=99, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
